$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- 1) After the paragraph containing the lone <w:tab/>, insert the
#        "Verificar librererías" paragraph (with spell-check proofing marks
#        around the misspelled word, exactly as the target XML has it). ---
$tabIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "`t`r") {
        $tabIndex = $i
    }
}

$tabPara = $d.Paragraphs.Item($tabIndex)
[void]$tabPara.Range.InsertParagraphAfter()
$verificarPara = $d.Paragraphs.Item($tabIndex + 1)
$verificarXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t xml:space="preserve">Verificar </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>librererías</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
[void]$verificarPara.Range.InsertXML($verificarXml)

# --- 2) Locate the last paragraph of the body (immediately before the
#        sectPr) and append the three new Odoo-12-dependency paragraphs
#        after it, each tagged en-US and carrying the proofing marks Word
#        would generate for the English shell commands. ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
[void]$lastPara.Range.InsertParagraphAfter()
$depsPara = $d.Paragraphs.Item($lastIndex + 1)

$cdXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>cd</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> /bin/lib</w:t></w:r>' +
    '</w:p>'

$aptXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>sudo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> apt-get install </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ttf</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>mscorefonts</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-installer</w:t></w:r>' +
    '</w:p>'

$emptyXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '</w:p>'

[void]$depsPara.Range.InsertXML($cdXml + $aptXml + $emptyXml)
